$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell value edits: language changed from French (fra) to English (eng),
#    descriptions translated, and the is_active flag re-authored as literal
#    text "TRUE" (was a boolean) for rows 2-5.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "eng"
$ws.Range("A3").Value = "eng"
$ws.Range("A4").Value = "eng"
$ws.Range("A5").Value = "eng"

$ws.Range("C2").Value = "Text File"
$ws.Range("C3").Value = "XML File"
$ws.Range("C4").Value = "Json File"
$ws.Range("C5").Value = "html file"

# Leading apostrophe forces these to be stored as literal text "TRUE"
# instead of being auto-coerced into a boolean TRUE value.
$ws.Range("D2").Value = "'TRUE"
$ws.Range("D3").Value = "'TRUE"
$ws.Range("D4").Value = "'TRUE"
$ws.Range("D5").Value = "'TRUE"

# ---------------------------------------------------------------------------
# 2. Header row (row 1) reformatting: simplified thin uniform border,
#    wrap text removed (alignment stays centered / top).
# ---------------------------------------------------------------------------
$header = $ws.Range("A1:D1")
$header.WrapText = $false
$header.Borders.Weight = 2
$header.Borders.Color = 0

# ---------------------------------------------------------------------------
# 3. Data rows (2-5): drop the explicit font/border/shading formatting that
#    used to distinguish them, reverting to the workbook's Normal style.
# ---------------------------------------------------------------------------
$ws.Range("A2:C5").Style = "Normal"
$ws.Range("D2:D5").Style = "Normal"

# Column D (including the header cell) is given a Text number format.
$ws.Range("D1:D5").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 4. Row heights: clear the explicit row heights / thick-bottom borders so
#    rows fall back to the sheet's default height.
# ---------------------------------------------------------------------------
$ws.Rows("1:5").EntireRow.AutoFit()

# ---------------------------------------------------------------------------
# 5. Sheet view / window state.
# ---------------------------------------------------------------------------
$ws.Range("M7").Select()
$excel.ActiveWindow.Zoom = 100

Write-Output "edit applied"
